$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.758.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.349.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.32"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.47"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +16.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.700.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.35"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.924"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.355.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.667.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.59"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.14"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.64"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.04"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.134"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0755"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.78"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.58"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.97%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0278"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.52"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +17.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.39"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +10.51%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.01"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.41%  "
